{"js": "// 1) Update the \"Sprint 6\" bold heading line: change the duration from\n//    \"1 Semana\" to \"2 Semanas\" and the date range from \"22-28 enero 2025\"\n//    to \"22 enero \u2013 04 febrero 2025\".\nconst sprint6 = context.document.body.search(\n  \"Gestionar dep\u00f3sitos en MXN. (1 Semana) 22-28 enero 2025\",\n  { matchCase: true }\n);\nsprint6.load(\"text\");\nawait context.sync();\n\nif (sprint6.items.length > 0) {\n  sprint6.items[0].insertText(\n    \"Gestionar dep\u00f3sitos en MXN. (2 Semanas) 22 enero \u2013 04 febrero 2025\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// 2) Remove the \"Registrar como transacci\u00f3n...\" bullet paragraph entirely.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  if (\n    paragraphs.items[i].text.indexOf(\n      \"Registrar como transacci\u00f3n en la tabla de transacciones\"\n    ) !== -1\n  ) {\n    paragraphs.items[i].delete();\n  }\n}\nawait context.sync();\n\n// 3) Drop the stale <w:lastRenderedPageBreak/> rendering marker on the\n//    \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \" run. Office.js has no\n//    direct property for this cached layout marker, so force a genuine\n//    text mutation (insert a sentinel, then restore the original text) on\n//    just that run's range so the run gets re-serialized without it.\nconst marker = context.document.body.search(\n  \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \",\n  { matchCase: true }\n);\nmarker.load(\"text\");\nawait context.sync();\n\nif (marker.items.length > 0) {\n  const target = marker.items[0];\n  target.insertText(\"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \\u0000\", \"Replace\");\n  await context.sync();\n\n  const resetSearch = context.document.body.search(\n    \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \\u0000\",\n    { matchCase: true }\n  );\n  await context.sync();\n  resetSearch.items[0].insertText(\n    \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the \"Sprint 6\" bold heading line: change the duration from\n#    \"1 Semana\" to \"2 Semanas\" and the date range from \"22-28 enero 2025\"\n#    to \"22 enero \u2013 04 febrero 2025\".\n$rng1 = $d.Content\n$rng1.Find.Text = \"Gestionar dep\u00f3sitos en MXN. (1 Semana) 22-28 enero 2025\"\n$rng1.Find.Replacement.Text = \"Gestionar dep\u00f3sitos en MXN. (2 Semanas) 22 enero \u2013 04 febrero 2025\"\n$rng1.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, \"wdReplaceAll\")\n\n# 2) Remove the \"Registrar como transacci\u00f3n...\" bullet paragraph entirely.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Registrar como transacci\u00f3n en la tabla de transacciones*\") {\n        $p.Range.Delete()\n    }\n}\n\n# 3) Drop the stale <w:lastRenderedPageBreak/> rendering marker sitting on\n#    the \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \" run. The COM object\n#    model has no property for this cached layout marker, so force a real\n#    text mutation (insert a sentinel, then restore the original text) on\n#    just that run's text so it gets re-serialized without the marker.\n$rng2 = $d.Content\n$rng2.Find.Text = \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \"\n$rng2.Find.Replacement.Text = \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el `u{0}\"\n$rng2.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, \"wdReplaceAll\")\n\n$rng3 = $d.Content\n$rng3.Find.Text = \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el `u{0}\"\n$rng3.Find.Replacement.Text = \"Se cre\u00f3 un m\u00e9todo as\u00edncrono activando el \"\n$rng3.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, \"wdReplaceAll\")\n"}
